$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D21: value changed from "administration" to "Coding/implementation" ---
$ws.Range("D21").Value = "Coding/implementation"

# --- Row 22 ---
$ws.Range("A22").Value = 45068
$ws.Range("B22").Formula = "=WEEKNUM(A22)"
$ws.Range("C22").Value = 3.25
$ws.Range("D22").Value = "Coding/implementation"
$ws.Range("E22").Value = "newsletter fini"
$ws.Range("F22").Value = "aucun problème"

# --- Row 23 ---
$ws.Range("A23").Value = 45068
$ws.Range("B23").Formula = "=WEEKNUM(A23)"
$ws.Range("C23").Value = 4.5
$ws.Range("D23").Value = "Coding/implementation"
$ws.Range("E23").Value = "contact fini"
$ws.Range("F23").Value = "aucun problème"

# --- Row 24 ---
$ws.Range("A24").Value = 45069
$ws.Range("B24").Formula = "=WEEKNUM(A24)"
$ws.Range("C24").Value = 2.25
$ws.Range("D24").Value = "Coding/implementation"
$ws.Range("E24").Value = "Gallery frontend fini"
$ws.Range("F24").Value = "aucun problème"

# --- Row 25 ---
$ws.Range("A25").Value = 45069
$ws.Range("B25").Formula = "=WEEKNUM(A25)"
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = "Coding/implementation"
$ws.Range("E25").Value = "Gallery backend fini"
$ws.Range("F25").Value = "aucun problème"

# --- Row 26 ---
$ws.Range("A26").Value = 45069
$ws.Range("B26").Formula = "=WEEKNUM(A26)"
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = "Coding/implementation"
$ws.Range("E26").Value = "Gallery modify fini"
$ws.Range("F26").Value = "aucun problème"

# --- Row 27 ---
$ws.Range("A27").Value = 45071
$ws.Range("B27").Formula = "=WEEKNUM(A27)"
$ws.Range("C27").Value = 5.25
$ws.Range("D27").Value = "Documentation"
$ws.Range("E27").Value = "documentation general"
$ws.Range("F27").Value = "aucun problème"

# --- Row 28 ---
$ws.Range("A28").Value = 45072
$ws.Range("B28").Formula = "=WEEKNUM(A28)"
$ws.Range("C28").Value = 5.25
$ws.Range("D28").Value = "Documentation"
$ws.Range("F28").Value = "Il n'y avait pas de diagramme avec les exigences que je voulais utiliser, j'ai donc utilisé UML mais cela ne correspond pas aux normes."
$ws.Range("E28").Value = "Structure de code de site"

# --- Update the Temps Total formula to cover the new rows ---
$ws.Range("I2").Formula = "=SUM(C2:C26)"

# --- Copy the formatting of row 21 (A:F) down into the new rows 22-28 ---
# (done after the values/formulas above so the dependent SUM formula keeps a
#  fresh cached result instead of the stale pre-paste snapshot)
$ws.Range("A21:F21").Copy()
$ws.Range("A22:F28").PasteSpecial(-4122)  # xlPasteFormats

# --- Update selected cell to match the final state ---
$ws.Range("C28").Select()
